$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 13 (they are removed entirely in the new version)
$ws.Range("A3:F13").EntireRow.Delete()

# Update row 2: only A2 (Comenzar) and D2 (Fin) keep values; B2/C2/E2/F2 are cleared
$ws.Range("A2").Value = "2024-10-25 02:38:09"
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = "2024-10-25 02:38:25"
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
